$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove one of the repeated empty (bold) paragraphs that sits
#    right before the "Arnau Marcos Almansa" right-aligned paragraph.
# ------------------------------------------------------------------
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $next = $d.Paragraphs.Item($i + 1)
    if ($p.Range.Text.TrimEnd([char]13) -eq "" -and
        $next.Range.Text.StartsWith("Arnau Marcos Almansa")) {
        $target = $p
        break
    }
}
$d.Range($target.Range.Start, $target.Range.End).Delete()

# ------------------------------------------------------------------
# 2) Locate the "Grup 45" paragraph and replace it together with the
#    following paragraph mark by two right-aligned paragraphs:
#      "Grup 45"      (unchanged text)
#      "Curs 2019-20" (new paragraph, added right below)
#    Replacing the whole original paragraph range (rather than just
#    inserting after it) lets the new paragraph marks start out with
#    no explicit run-properties, matching a freshly typed line.
# ------------------------------------------------------------------
$grup = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.TrimEnd([char]13) -eq "Grup 45") {
        $grup = $p
        break
    }
}

$rng = $d.Range($grup.Range.Start, $grup.Range.End)
$xml = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="right"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial" w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="24"/><w:lang w:val="ca-ES"/></w:rPr><w:t>Grup 45</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Normal"/><w:jc w:val="right"/><w:rPr/></w:pPr><w:r><w:rPr><w:rFonts w:cs="Arial" w:ascii="Arial" w:hAnsi="Arial"/><w:sz w:val="24"/><w:lang w:val="ca-ES"/></w:rPr><w:t>Curs 2019-20</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$null = $rng.InsertXML($xml)
